$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Bidang" header in column I
$ws.Range("I1").Value = "Bidang"

# Match the header formatting used by the rest of row 1 (fill/border/centered)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Set the width for the new column I (closest attainable to the template's 13.6328125)
$ws.Columns.Item(9).ColumnWidth = 12.8

# Update the active selection to match the target state
$ws.Range("F13").Select()
